$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 to make room for the group headers, pushing
# the existing header row (Name/IP Address/Username/Password/Port) to row 2
$ws.Rows("1:1").Insert()

# --- Row 1: group headers (values) ---
$ws.Range("A1").Value = "Common Data"
$ws.Range("G1").Value = "Autoprovision"
$ws.Range("K1").Value = "Edit Wave"
$ws.Range("N1").Value = "NIC's"
$ws.Range("V1").Value = "Edit Wave"

# --- Row 2: detail column headers ---
$ws.Range("F2").Value = "Wave Name"
$ws.Range("G2").Value = "Environment"
$ws.Range("H2").Value = "Cluster Name"
$ws.Range("I2").Value = "ESX Host"
$ws.Range("J2").Value = "Datastore"
$ws.Range("K2").Value = "Cluster Name"
$ws.Range("L2").Value = "ESX Host"
$ws.Range("M2").Value = "Datastore"
$ws.Range("N2").Value = "Device Name"
$ws.Range("O2").Value = "Type"
$ws.Range("P2").Value = "Network Name"
$ws.Range("Q2").Value = "DHCP / Static IP"
$ws.Range("R2").Value = "CIDR"
$ws.Range("S2").Value = "Gateway"
$ws.Range("T2").Value = "DNS1"
$ws.Range("U2").Value = "DNS2"
$ws.Range("V2").Value = "VM Folder"
$ws.Range("W2").Value = "Resource Pool"
$ws.Range("X2").Value = "Routes"

# --- Format + merge the uniform (fully bold+centered) group-header blocks ---
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").Merge()

$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

$ws.Range("G1:J1").Font.Bold = $true
$ws.Range("G1:J1").HorizontalAlignment = -4108
$ws.Range("G1:J1").Merge()

$ws.Range("K1:M1").Font.Bold = $true
$ws.Range("K1:M1").HorizontalAlignment = -4108
$ws.Range("K1:M1").Merge()

$ws.Range("N1:U1").Font.Bold = $true
$ws.Range("N1:U1").HorizontalAlignment = -4108
$ws.Range("N1:U1").Merge()

# --- "Edit Wave" (V1:X1): merge first, then format the anchor, so the two
# trailing cells (W1/X1) can be given a different, non-bold centered style ---
$ws.Range("V1:X1").Merge()
$ws.Range("V1").Font.Bold = $true
$ws.Range("V1").HorizontalAlignment = -4108
$ws.Cells.Item(1, 23).HorizontalAlignment = -4108
$ws.Cells.Item(1, 24).HorizontalAlignment = -4108

# --- Column widths to fit the new content ---
$ws.Range("A1:X2").Columns.AutoFit()

# --- Selection matches the post-edit state ---
$ws.Range("F11").Select()
